$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Find the first empty row below the existing GSC export data (column A).
$lastCell = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162)
$newRow = $lastCell.Row + 1

# Column A holds the date as literal text (not a real Excel date), matching
# the rest of the export. Writing the text directly via .Value/.Value2 gets
# auto-coerced into a date serial by Excel's locale-aware input parsing, so
# route it through a text formula and flatten it back to a static value
# (copy / paste-special-values) to keep it a plain shared string with the
# workbook's default style.
$ws.Cells.Item($newRow, 1).Formula = '="2025-11-29"'
$ws.Cells.Item($newRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163)

$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 27
